$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 4191701.42
$ws.Range("C9").Value = 659076.24
$ws.Range("D9").Value = 4850777.66
$ws.Range("E9").Value = 13.58702225077865
$ws.Range("F9").Value = 86.41297774922134
$ws.Range("G9").Value = -36.30356770963679
$ws.Range("H9").Value = -24.30366629055501
$ws.Range("I9").Value = 42177
$ws.Range("J9").Value = 1814
$ws.Range("K9").Value = 43991
$ws.Range("L9").Value = 30540
$ws.Range("M9").Value = 158.8335841519319
$ws.Range("N9").Value = 8.438922979987096
